# 1. improved comments 2. Stopped browser from opening 3. made error handling more clear
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "invalid user" result is now reported as its own leaderboard row,
# right under the header, instead of silently being skipped. Insert a row
# so the rest of the leaderboard shifts down instead of being overwritten.
$ws.Rows("3").Insert()

$ws.Range("A3").Value = "['Nithin_XS1223']"
$ws.Range("B3").Value = "INVALID USER"

# The rating column used to store each contestant's rating as text; store
# real numbers now so the values can be sorted/compared correctly.
$ratings = @(3759, 3697, 3662, 3644, 3505, 3486)
for ($i = 0; $i -lt $ratings.Length; $i++) {
    $row = 4 + $i
    $ws.Range("B$row").Value = $ratings[$i]
}
$ws.Range("B4:B9").NumberFormat = "General"

# Let the columns re-fit their contents now that the data has changed.
$ws.Columns("A:B").AutoFit()

# Leave the selection where the user last left off while testing.
$null = $ws.Range("I8").Select()
